$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45265, 16, 100114001, 'Papa', 'Asterix', '1a nueva(o)', 200, 23000, 23000, 23000, '$/saco 25 kilos', 'Región del Maule', 920, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45265, 16, 100114001, 'Papa', 'Asterix', '2a nueva(o)', 200, 20000, 20000, 20000, '$/saco 25 kilos', 'Región del Maule', 800, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45222, 16, 100114001, 'Papa', 'Asterix', '1a (guarda lavada)', 200, 32000, 32000, 32000, '$/malla 25 kilos', 'Región de Los Lagos', 1280, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45222, 16, 100114001, 'Papa', 'Asterix', '1a (guarda)', 300, 30000, 30000, 30000, '$/saco 25 kilos', 'Región de Los Lagos', 1200, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44592, 16, 100114001, 'Papa', 'Patagonia', '1a nueva(o)', 120, 6500, 7000, 6750, '$/saco 25 kilos', 'Provincia de Diguillín', 270, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45128, 16, 100114001, 'Papa', 'Asterix', '1a (guarda lavada)', 100, 19000, 19000, 19000, '$/malla 25 kilos', 'Región de Los Lagos', 760, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45128, 16, 100114001, 'Papa', 'Asterix', '1a (guarda)', 150, 18000, 18000, 18000, '$/saco 25 kilos', 'Región de Los Lagos', 720, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45128, 16, 100114001, 'Papa', 'Asterix', '2a (guarda)', 150, 16000, 16000, 16000, '$/saco 25 kilos', 'Región de Los Lagos', 640, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44761, 16, 100114001, 'Papa', 'Asterix', '1a (guarda)', 160, 7000, 7500, 7250, '$/saco 25 kilos', 'Provincia de Diguillín', 290, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44761, 16, 100114001, 'Papa', 'Patagonia', '1a (guarda)', 120, 7500, 8000, 7750, '$/saco 25 kilos', 'Provincia de Diguillín', 310, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44363, 16, 100114001, 'Papa', 'Patagonia', '1a (cosecha)', 120, 5500, 6000, 5750, '$/saco 25 kilos', 'Provincia de Diguillín', 230, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45097, 16, 100114001, 'Papa', 'Asterix', '1a (guarda lavada)', 130, 13000, 14000, 13538, '$/malla 25 kilos', 'Región de Los Lagos', 542, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45097, 16, 100114001, 'Papa', 'Patagonia', '1a (guarda)', 200, 11000, 12000, 11500, '$/saco 25 kilos', 'Región de Los Lagos', 460, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44222, 16, 100114001, 'Papa', 'Asterix', '1a nueva(o)', 300, 8500, 9000, 8733, '$/saco 25 kilos', 'Región del Maule', 349, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44271, 16, 100114001, 'Papa', 'Rosara', '1a (cosecha)', 280, 6500, 7000, 6750, '$/saco 25 kilos', 'Región de Los Lagos', 270, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44880, 16, 100114001, 'Papa', 'Rosara', '1a nueva(o)', 120, 11000, 12000, 11500, '$/saco 25 kilos', 'Región de O''Higgins', 460, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44650, 16, 100114001, 'Papa', 'Patagonia', '1a (cosecha)', 200, 7000, 7500, 7250, '$/saco 25 kilos', 'Provincia de Diguillín', 290, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44399, 16, 100114001, 'Papa', 'Patagonia', '1a (cosecha)', 120, 6000, 6500, 6250, '$/saco 25 kilos', 'Provincia de Diguillín', 250, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45107, 16, 100114001, 'Papa', 'Asterix', '1a (guarda lavada)', 250, 17000, 18000, 17600, '$/malla 25 kilos', 'Región de Los Lagos', 704, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45107, 16, 100114001, 'Papa', 'Asterix', '1a (guarda)', 180, 16000, 16000, 16000, '$/saco 25 kilos', 'Región de Los Lagos', 640, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45107, 16, 100114001, 'Papa', 'Patagonia', '1a (guarda)', 180, 17000, 17000, 17000, '$/saco 25 kilos', 'Región de Los Lagos', 680, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44435, 16, 100114001, 'Papa', 'Patagonia', '1a (guarda lavada)', 600, 7500, 8000, 7750, '$/saco 25 kilos', 'Región de La Araucanía', 310, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44435, 16, 100114001, 'Papa', 'Patagonia', '1a (guarda)', 460, 6000, 7500, 6902, '$/saco 25 kilos', 'Provincia de Diguillín', 276, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44435, 16, 100114001, 'Papa', 'Rodeo', '1a (guarda)', 760, 6000, 7000, 6447, '$/saco 25 kilos', 'Provincia de Diguillín', 258, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44609, 16, 100114001, 'Papa', 'Patagonia', '1a nueva(o)', 200, 6500, 7000, 6750, '$/saco 25 kilos', 'Provincia de Diguillín', 270, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44509, 16, 100114001, 'Papa', 'Patagonia', '1a (guarda)', 360, 7000, 8000, 7500, '$/saco 25 kilos', 'Provincia de Diguillín', 300, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45050, 16, 100114001, 'Papa', 'Patagonia', '1a (cosecha lavada)', 150, 12000, 12000, 12000, '$/malla 25 kilos', 'Región de Los Lagos', 480, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45050, 16, 100114001, 'Papa', 'Patagonia', '1a (cosecha)', 600, 10000, 11000, 10500, '$/saco 25 kilos', 'Región de Los Lagos', 420, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45180, 16, 100114001, 'Papa', 'Asterix', '1a (guarda lavada)', 100, 32000, 32000, 32000, '$/malla 25 kilos', 'Región de Los Lagos', 1280, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45180, 16, 100114001, 'Papa', 'Asterix', '1a (guarda)', 200, 30000, 30000, 30000, '$/saco 25 kilos', 'Región de Los Lagos', 1200, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44824, 16, 100114001, 'Papa', 'Patagonia', '1a (guarda)', 120, 7500, 8000, 7750, '$/saco 25 kilos', 'Provincia de Diguillín', 310, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44547, 16, 100114001, 'Papa', 'Asterix', '1a nueva(o)', 1200, 8500, 9000, 8750, '$/saco 25 kilos', 'Región del Maule', 350, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44539, 16, 100114001, 'Papa', 'Asterix', '1a nueva(o)', 160, 9000, 9500, 9250, '$/saco 25 kilos', 'Región del Maule', 370, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44539, 16, 100114001, 'Papa', 'Patagonia', '1a nueva(o)', 160, 9000, 9500, 9250, '$/saco 25 kilos', 'Región del Maule', 370, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44894, 16, 100114001, 'Papa', 'Pukará', '1a nueva(o)', 120, 11000, 12000, 11500, '$/saco 25 kilos', 'Región del Maule', 460, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44894, 16, 100114001, 'Papa', 'Rosara', '1a nueva(o)', 120, 12000, 13000, 12500, '$/saco 25 kilos', 'Región de O''Higgins', 500, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44193, 16, 100114001, 'Papa', 'Asterix', '1a nueva(o)', 120, 13000, 14000, 13500, '$/saco 25 kilos', 'Provincia de Diguillín', 540, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44631, 16, 100114001, 'Papa', 'Asterix', '1a (cosecha lavada)', 100, 9000, 9000, 9000, '$/malla 25 kilos', 'Región Metropolitana', 360, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44631, 16, 100114001, 'Papa', 'Patagonia', '1a (cosecha)', 120, 6000, 6500, 6250, '$/saco 25 kilos', 'Provincia de Diguillín', 250, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44914, 16, 100114001, 'Papa', 'Asterix', '1a nueva(o)', 300, 12000, 12000, 12000, '$/saco 25 kilos', 'Región del Maule', 480, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44239, 16, 100114001, 'Papa', 'Patagonia', '1a nueva(o)', 120, 7000, 7500, 7250, '$/malla 25 kilos', 'Provincia de Diguillín', 290, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45168, 16, 100114001, 'Papa', 'Asterix', '1a (guarda lavada)', 250, 33000, 33000, 33000, '$/malla 25 kilos', 'Región de La Araucanía', 1320, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45168, 16, 100114001, 'Papa', 'Asterix', '1a (guarda)', 200, 32000, 32000, 32000, '$/saco 25 kilos', 'Región de Los Lagos', 1280, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44574, 16, 100114001, 'Papa', 'Asterix', '1a nueva(o)', 400, 7500, 8000, 7750, '$/saco 25 kilos', 'Región del Maule', 310, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44574, 16, 100114001, 'Papa', 'Asterix', '2a nueva(o)', 160, 6500, 7000, 6750, '$/saco 25 kilos', 'Región del Maule', 270, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44214, 16, 100114001, 'Papa', 'Asterix', '1a nueva(o)', 270, 8500, 9000, 8778, '$/saco 25 kilos', 'Región del Maule', 351, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44214, 16, 100114001, 'Papa', 'Rosara', '1a nueva(o)', 160, 8000, 8500, 8281, '$/saco 25 kilos', 'Región del Maule', 331, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44895, 16, 100114001, 'Papa', 'Pukará', '1a nueva(o)', 120, 11000, 12000, 11500, '$/saco 25 kilos', 'Región del Maule', 460, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44895, 16, 100114001, 'Papa', 'Rosara', '1a nueva(o)', 120, 12000, 13000, 12500, '$/saco 25 kilos', 'Región de O''Higgins', 500, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45135, 16, 100114001, 'Papa', 'Asterix', '1a (guarda lavada)', 100, 19000, 19000, 19000, '$/malla 25 kilos', 'Región de Los Lagos', 760, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45135, 16, 100114001, 'Papa', 'Asterix', '1a (guarda)', 100, 18000, 18000, 18000, '$/saco 25 kilos', 'Región de Los Lagos', 720, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44567, 16, 100114001, 'Papa', 'Asterix', '1a nueva(o)', 300, 7500, 8000, 7750, '$/saco 25 kilos', 'Región del Maule', 310, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44517, 16, 100114001, 'Papa', 'Patagonia', '1a nueva(o)', 60, 11000, 12000, 11500, '$/saco 25 kilos', 'Región del Maule', 460, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45121, 16, 100114001, 'Papa', 'Asterix', '1a (guarda lavada)', 100, 19000, 19000, 19000, '$/malla 25 kilos', 'Región de Los Lagos', 760, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45121, 16, 100114001, 'Papa', 'Asterix', '1a (guarda)', 150, 18000, 18000, 18000, '$/saco 25 kilos', 'Región de Los Lagos', 720, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44203, 16, 100114001, 'Papa', 'Asterix', '1a nueva(o)', 200, 11000, 12000, 11600, '$/saco 25 kilos', 'Región del Maule', 464, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44818, 16, 100114001, 'Papa', 'Patagonia', '1a (guarda)', 120, 7000, 7500, 7250, '$/saco 25 kilos', 'Región de Ñuble', 290, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44293, 16, 100114001, 'Papa', 'Patagonia', '1a (cosecha)', 120, 6500, 7000, 6750, '$/saco 25 kilos', 'Provincia de Diguillín', 270, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 44189, 16, 100114001, 'Papa', 'Asterix', '1a nueva(o)', 120, 14000, 15000, 14500, '$/saco 25 kilos', 'Región del Maule', 580, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45100, 16, 100114001, 'Papa', 'Asterix', '1a (guarda lavada)', 100, 13000, 13000, 13000, '$/malla 25 kilos', 'Región de Los Lagos', 520, 25, 'Hortaliza')
    ,@(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45100, 16, 100114001, 'Papa', 'Patagonia', '1a (guarda)', 100, 12000, 12000, 12000, '$/saco 25 kilos', 'Región de Los Lagos', 480, 25, 'Hortaliza')
)

$dateFormat = $ws.Range("D792").NumberFormat

$startRow = 734
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $row = $newData[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $row[$c]
    }
    if ($r -gt 792) {
        $ws.Cells.Item($r, 4).NumberFormat = $dateFormat
    }
}

Write-Host "Rows written: $($newData.Count)"